# Apply the commit's changes:
#  - move the selection on the first sheet (Tabelle1) to D15
#  - add a new worksheet "Sheet4" at the end of the workbook with a small
#    Col1/Col2/Col3 table, and leave it as the active sheet/tab with F13
#    selected (which also bumps the workbook's activeTab to 3).

$wb = $excel.ActiveWorkbook

# --- Tabelle1: move the selection from E12 to D15, no longer the active tab ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("D15").Select()

# --- Add "Sheet4" as the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Sheet4"

# --- Populate the new sheet ---
$newSheet.Range("A1").Value = "Col1"
$newSheet.Range("B1").Value = "Col2"
$newSheet.Range("C1").Value = "Col3"
$newSheet.Range("A2").Value = "text2"
$newSheet.Range("B2").Value = "text1"
$newSheet.Range("C2").Value = "text3"

# --- Make Sheet4 the active sheet/tab with F13 selected ---
$newSheet.Range("F13").Select()
